# Update performance document (sheet "20121215D"):
#  - Header row: replace the stale "Optimized type classifications / id /
#    read-only string optimization / mangling optimization" header set with
#    the two current milestones ("Read only string in parser and symbol" is
#    already in C1 via shared string reuse; D1 gets the new label) and blank
#    out the now-unused E1:G1 headers.
#  - Fill in the "Read only string in parser and symbol" run's raw data
#    (column D, rows 2-11) which previously had no samples (hence the
#    #DIV/0! errors in the derived rows 12-16).
#  - Move the active selection to D15 (last place edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("20121215D")

# Header row: C1 keeps its text ("Read only string in parser and symbol"),
# D1 gets the new run label, E1:G1 are cleared (no more runs recorded yet).
$ws.Range("D1").Value = "Remove unused symbol insertion when parameter is anonymous"
$ws.Range("E1").Value = ""
$ws.Range("F1").Value = ""
$ws.Range("G1").Value = ""

# Raw per-iteration timings for the "Read only string in parser and symbol"
# column (D2:D11). Filling these in lets the existing shared AVERAGE/VAR.S/
# T.TEST/ratio formulas in rows 12-16 recalculate instead of erroring out.
$ws.Range("D2").Value = 4765
$ws.Range("D3").Value = 4747
$ws.Range("D4").Value = 4771
$ws.Range("D5").Value = 4754
$ws.Range("D6").Value = 4751
$ws.Range("D7").Value = 4812
$ws.Range("D8").Value = 4771
$ws.Range("D9").Value = 4795
$ws.Range("D10").Value = 4736
$ws.Range("D11").Value = 4780

# Leave the selection where the author last left it.
$ws.Activate()
$ws.Range("D15").Select()
